# Insert a new weekly data row at row 73 (pushing existing rows 73:200 down to 74:201),
# matching the commit "Fruta / hortaliza, semanal" which adds one new weekly price
# observation into the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73:200 down by inserting a fresh blank row at position 73.
$ws.Rows("73:73").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value = "12/03/2021"
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112023
$ws.Range("G73").Value = "Brócoli"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 300
$ws.Range("K73").Value = 700
$ws.Range("L73").Value = 800
$ws.Range("M73").Value = 750
$ws.Range("N73").Value = "$/unidad"
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 750
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"
